$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.602.50'
$ws.Range("E2").Value = '  +5.76%  '

$ws.Range("D3").Value = '2.744.00'
$ws.Range("E3").Value = '  +4.05%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '116.73'
$ws.Range("E5").Value = '  +5.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '332.60'
$ws.Range("E6").Value = '  +3.07%  '

$ws.Range("E7").Value = '  +2.42%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.575'
$ws.Range("E9").Value = '  +6.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.54'
$ws.Range("E10").Value = '  +5.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.12'
$ws.Range("E11").Value = '  +1.16%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0831'
$ws.Range("E12").Value = '  +2.60%  '

$ws.Range("E13").Value = '  +2.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.57'
$ws.Range("E14").Value = '  +4.82%  '

$ws.Range("D15").Value = '3.179.55'
$ws.Range("E15").Value = '  +4.42%  '

$ws.Range("D16").Value = '2.746.17'
$ws.Range("E16").Value = '  +3.77%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.884'
$ws.Range("E17").Value = '  +2.31%  '

$ws.Range("D18").Value = '51.548.00'
$ws.Range("E18").Value = '  +5.66%  '

$ws.Range("B19").Value = 'ImmutableX'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.09'
$ws.Range("E19").Value = '  +7.17%  '

$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.44'
$ws.Range("E20").Value = '  +4.58%  '

$ws.Range("E21").Value = '  +2.00%  '

$ws.Range("D22").Value = '0.0₃0961'
$ws.Range("E22").Value = '  +2.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '277.98'
$ws.Range("E23").Value = '  +2.97%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.54'
$ws.Range("E24").Value = '  +1.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.64'
$ws.Range("E25").Value = '  +4.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.70'
$ws.Range("E26").Value = '  +2.32%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.21'
$ws.Range("E28").Value = '  +0.95%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.22'
$ws.Range("E29").Value = '  -0.24%  '

$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.139'
$ws.Range("E30").Value = '  +1.72%  '

$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.11'
$ws.Range("E31").Value = '  +0.32%  '

$ws.Range("B32").Value = 'OKB'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.27'
$ws.Range("E32").Value = '  +1.56%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.56'
$ws.Range("E33").Value = '  +1.55%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0818'
$ws.Range("E34").Value = '  +2.74%  '

$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.07'
$ws.Range("E36").Value = '  -1.19%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.08'
$ws.Range("E37").Value = '  +2.21%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.93'
$ws.Range("E38").Value = '  -0.10%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.16'
$ws.Range("E39").Value = '  -0.18%  '

$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '129.50'
$ws.Range("E40").Value = '  +3.54%  '

$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.23'
$ws.Range("E41").Value = '  +2.62%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0345'
$ws.Range("E42").Value = '  +9.58%  '

$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.113'
$ws.Range("E43").Value = '  +2.42%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.26'
$ws.Range("E44").Value = '  +5.06%  '

$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.41'
$ws.Range("E45").Value = '  +13.36%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.106.93'
$ws.Range("E46").Value = '  +1.66%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.33'
$ws.Range("E47").Value = '  +3.13%  '

$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.23'
$ws.Range("E48").Value = '  +2.29%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.55'
$ws.Range("E49").Value = '  +7.57%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.94'
$ws.Range("E50").Value = '  -0.46%  '

$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '60.03'
$ws.Range("E51").Value = '  +2.31%  '
